# Auto-generated edit script: update cached market-price derived values
# in the per-job "Profits" worksheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1828.2778  # H17
$ws.Cells.Item(17, 10).Value = 1828.2778  # J17
$ws.Cells.Item(17, 12).Value = 5484.8334  # L17
$ws.Cells.Item(17, 14).Value = -5820.8334  # N17
$ws.Cells.Item(19, 8).Value = 3325  # H19
$ws.Cells.Item(19, 9).Value = 3898  # I19
$ws.Cells.Item(19, 10).Value = 2943  # J19
$ws.Cells.Item(19, 11).Value = 3898  # K19
$ws.Cells.Item(19, 12).Value = 2943  # L19
$ws.Cells.Item(19, 13).Value = -3723  # M19
$ws.Cells.Item(19, 14).Value = -3293  # N19
$ws.Cells.Item(28, 8).Value = 1858.5714  # H28
$ws.Cells.Item(28, 9).Value = 1858.5714  # I28
$ws.Cells.Item(28, 11).Value = 1858.5714  # K28
$ws.Cells.Item(28, 13).Value = -1373.5714  # M28
$ws.Cells.Item(40, 8).Value = 1720.5  # H40
$ws.Cells.Item(40, 9).Value = 1653.4546  # I40
$ws.Cells.Item(40, 11).Value = 1653.4546  # K40
$ws.Cells.Item(40, 13).Value = -1478.4546  # M40
$ws.Cells.Item(43, 8).Value = 3121.3333  # H43
$ws.Cells.Item(43, 9).Value = 1632  # I43
$ws.Cells.Item(43, 11).Value = 1632  # K43
$ws.Cells.Item(43, 13).Value = -1563  # M43
$ws.Cells.Item(53, 8).Value = 259.05884  # H53
$ws.Cells.Item(53, 9).Value = 176.11111  # I53
$ws.Cells.Item(53, 11).Value = 176.11111  # K53
$ws.Cells.Item(53, 13).Value = 460.88889  # M53
$ws.Cells.Item(70, 8).Value = 2324.75  # H70
$ws.Cells.Item(70, 9).Value = 1300  # I70
$ws.Cells.Item(70, 10).Value = 2666.3333  # J70
$ws.Cells.Item(70, 11).Value = 3900  # K70
$ws.Cells.Item(70, 12).Value = 7998.999899999999  # L70
$ws.Cells.Item(70, 13).Value = -3630  # M70
$ws.Cells.Item(70, 14).Value = -8538.999899999999  # N70
$ws.Cells.Item(73, 8).Value = 2324.75  # H73
$ws.Cells.Item(73, 9).Value = 1300  # I73
$ws.Cells.Item(73, 10).Value = 2666.3333  # J73
$ws.Cells.Item(73, 11).Value = 3900  # K73
$ws.Cells.Item(73, 12).Value = 7998.999899999999  # L73
$ws.Cells.Item(73, 13).Value = -2964  # M73
$ws.Cells.Item(73, 14).Value = -9870.999899999999  # N73
$ws.Cells.Item(76, 8).Value = 3699  # H76
$ws.Cells.Item(76, 9).Value = 3749.1667  # I76
$ws.Cells.Item(76, 10).Value = 3398  # J76
$ws.Cells.Item(76, 11).Value = 3749.1667  # K76
$ws.Cells.Item(76, 12).Value = 3398  # L76
$ws.Cells.Item(76, 13).Value = -3434.1667  # M76
$ws.Cells.Item(76, 14).Value = -4028  # N76
$ws.Cells.Item(79, 8).Value = 3699  # H79
$ws.Cells.Item(79, 9).Value = 3749.1667  # I79
$ws.Cells.Item(79, 10).Value = 3398  # J79
$ws.Cells.Item(79, 11).Value = 3749.1667  # K79
$ws.Cells.Item(79, 12).Value = 3398  # L79
$ws.Cells.Item(79, 13).Value = -2657.1667  # M79
$ws.Cells.Item(79, 14).Value = -5582  # N79
$ws.Cells.Item(80, 8).Value = 5925.25  # H80
$ws.Cells.Item(80, 9).Value = 5625.5  # I80
$ws.Cells.Item(80, 10).Value = 6225  # J80
$ws.Cells.Item(80, 11).Value = 16876.5  # K80
$ws.Cells.Item(80, 12).Value = 18675  # L80
$ws.Cells.Item(80, 13).Value = -15878.5  # M80
$ws.Cells.Item(80, 14).Value = -20671  # N80
$ws.Cells.Item(83, 8).Value = 5925.25  # H83
$ws.Cells.Item(83, 9).Value = 5625.5  # I83
$ws.Cells.Item(83, 10).Value = 6225  # J83
$ws.Cells.Item(83, 11).Value = 50629.5  # K83
$ws.Cells.Item(83, 12).Value = 56025  # L83
$ws.Cells.Item(83, 13).Value = -45637.5  # M83
$ws.Cells.Item(83, 14).Value = -66009  # N83
$ws.Cells.Item(88, 8).Value = 2200.25  # H88
$ws.Cells.Item(88, 10).Value = 2200.25  # J88
$ws.Cells.Item(88, 12).Value = 2200.25  # L88
$ws.Cells.Item(88, 14).Value = -3012.25  # N88
$ws.Cells.Item(91, 8).Value = 2200.25  # H91
$ws.Cells.Item(91, 10).Value = 2200.25  # J91
$ws.Cells.Item(91, 12).Value = 2200.25  # L91
$ws.Cells.Item(91, 14).Value = -5008.25  # N91
$ws.Cells.Item(92, 8).Value = 238.16667  # H92
$ws.Cells.Item(92, 9).Value = 238.16667  # I92
$ws.Cells.Item(92, 10).Value = 0  # J92
$ws.Cells.Item(92, 11).Value = 238.16667  # K92
$ws.Cells.Item(92, 12).Value = 0  # L92
$ws.Cells.Item(92, 13).Value = 1009.83333  # M92
$ws.Cells.Item(92, 14).ClearContents()  # N92 removed
$ws.Cells.Item(98, 8).Value = 1277.0769  # H98
$ws.Cells.Item(98, 9).Value = 955.7778  # I98
$ws.Cells.Item(98, 10).Value = 2000  # J98
$ws.Cells.Item(98, 11).Value = 955.7778  # K98
$ws.Cells.Item(98, 12).Value = 2000  # L98
$ws.Cells.Item(98, 13).Value = 542.2222  # M98
$ws.Cells.Item(98, 14).Value = -4996  # N98
$ws.Cells.Item(103, 8).Value = 942  # H103
$ws.Cells.Item(103, 9).Value = 749.5  # I103
$ws.Cells.Item(103, 10).Value = 1019  # J103
$ws.Cells.Item(103, 11).Value = 2248.5  # K103
$ws.Cells.Item(103, 12).Value = 3057  # L103
$ws.Cells.Item(103, 13).Value = -1662.5  # M103
$ws.Cells.Item(103, 14).Value = -4229  # N103
$ws.Cells.Item(106, 8).Value = 36438.2  # H106
$ws.Cells.Item(106, 9).Value = 36438.2  # I106
$ws.Cells.Item(106, 11).Value = 36438.2  # K106
$ws.Cells.Item(106, 13).Value = -35807.2  # M106
$ws.Cells.Item(107, 8).Value = 692.8570999999999  # H107
$ws.Cells.Item(107, 9).Value = 661.6667  # I107
$ws.Cells.Item(107, 11).Value = 661.6667  # K107
$ws.Cells.Item(107, 13).Value = 1258.3333  # M107
$ws.Cells.Item(113, 8).Value = 2270.8  # H113
$ws.Cells.Item(113, 9).Value = 1974.5  # I113
$ws.Cells.Item(113, 10).Value = 2468.3333  # J113
$ws.Cells.Item(113, 11).Value = 1974.5  # K113
$ws.Cells.Item(113, 12).Value = 2468.3333  # L113
$ws.Cells.Item(113, 13).Value = 1279.5  # M113
$ws.Cells.Item(113, 14).Value = -8976.3333  # N113
$ws.Cells.Item(116, 8).Value = 1799  # H116
$ws.Cells.Item(116, 10).Value = 1799  # J116
$ws.Cells.Item(116, 12).Value = 1799  # L116
$ws.Cells.Item(116, 14).Value = -8683  # N116
$ws.Cells.Item(118, 8).Value = 896  # H118
$ws.Cells.Item(118, 9).Value = 896  # I118
$ws.Cells.Item(118, 11).Value = 2688  # K118
$ws.Cells.Item(118, 13).Value = -1031  # M118
$ws.Cells.Item(122, 8).Value = 1277.0769  # H122
$ws.Cells.Item(122, 9).Value = 955.7778  # I122
$ws.Cells.Item(122, 10).Value = 2000  # J122
$ws.Cells.Item(122, 11).Value = 2867.3334  # K122
$ws.Cells.Item(122, 12).Value = 6000  # L122
$ws.Cells.Item(122, 13).Value = -417.3334  # M122
$ws.Cells.Item(122, 14).Value = -10900  # N122
$ws.Cells.Item(127, 8).Value = 2098.5  # H127
$ws.Cells.Item(127, 10).Value = 1886  # J127
$ws.Cells.Item(127, 12).Value = 5658  # L127
$ws.Cells.Item(127, 14).Value = -15578  # N127
$ws.Cells.Item(132, 8).Value = 1004.375  # H132
$ws.Cells.Item(132, 9).Value = 906.1667  # I132
$ws.Cells.Item(132, 11).Value = 2718.5001  # K132
$ws.Cells.Item(132, 13).Value = -188.5001000000002  # M132
$ws.Cells.Item(138, 8).Value = 3676.9302  # H138
$ws.Cells.Item(138, 9).Value = 1719.4445  # I138
$ws.Cells.Item(138, 11).Value = 5158.333500000001  # K138
$ws.Cells.Item(138, 13).Value = -18.33350000000064  # M138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2141375.8  # H32
$ws.Cells.Item(32, 9).Value = 2502608.8  # I32
$ws.Cells.Item(32, 11).Value = 2502608.8  # K32
$ws.Cells.Item(32, 13).Value = -2502321.8  # M32
$ws.Cells.Item(45, 8).Value = 144026.14  # H45
$ws.Cells.Item(45, 9).Value = 251247.5  # I45
$ws.Cells.Item(45, 10).Value = 1064.3334  # J45
$ws.Cells.Item(45, 11).Value = 251247.5  # K45
$ws.Cells.Item(45, 12).Value = 1064.3334  # L45
$ws.Cells.Item(45, 13).Value = -250870.5  # M45
$ws.Cells.Item(45, 14).Value = -1818.3334  # N45
$ws.Cells.Item(56, 8).Value = 11921.143  # H56
$ws.Cells.Item(56, 9).Value = 8916.666999999999  # I56
$ws.Cells.Item(56, 10).Value = 29948  # J56
$ws.Cells.Item(56, 11).Value = 8916.666999999999  # K56
$ws.Cells.Item(56, 12).Value = 29948  # L56
$ws.Cells.Item(56, 13).Value = -8174.666999999999  # M56
$ws.Cells.Item(56, 14).Value = -31432  # N56
$ws.Cells.Item(61, 8).Value = 2095.889  # H61
$ws.Cells.Item(61, 9).Value = 2188.3125  # I61
$ws.Cells.Item(61, 11).Value = 2188.3125  # K61
$ws.Cells.Item(61, 13).Value = -1976.3125  # M61
$ws.Cells.Item(63, 8).Value = 5089.5713  # H63
$ws.Cells.Item(63, 9).Value = 4988.769  # I63
$ws.Cells.Item(63, 11).Value = 4988.769  # K63
$ws.Cells.Item(63, 13).Value = -4302.769  # M63
$ws.Cells.Item(66, 8).Value = 5089.5713  # H66
$ws.Cells.Item(66, 9).Value = 4988.769  # I66
$ws.Cells.Item(66, 11).Value = 24943.845  # K66
$ws.Cells.Item(66, 13).Value = -21511.845  # M66
$ws.Cells.Item(88, 8).Value = 3415.5386  # H88
$ws.Cells.Item(88, 9).Value = 3500.8  # I88
$ws.Cells.Item(88, 10).Value = 3362.25  # J88
$ws.Cells.Item(88, 11).Value = 3500.8  # K88
$ws.Cells.Item(88, 12).Value = 3362.25  # L88
$ws.Cells.Item(88, 13).Value = -3094.8  # M88
$ws.Cells.Item(88, 14).Value = -4174.25  # N88
$ws.Cells.Item(91, 8).Value = 3415.5386  # H91
$ws.Cells.Item(91, 9).Value = 3500.8  # I91
$ws.Cells.Item(91, 10).Value = 3362.25  # J91
$ws.Cells.Item(91, 11).Value = 3500.8  # K91
$ws.Cells.Item(91, 12).Value = 3362.25  # L91
$ws.Cells.Item(91, 13).Value = -2096.8  # M91
$ws.Cells.Item(91, 14).Value = -6170.25  # N91
$ws.Cells.Item(136, 8).Value = 2095.889  # H136
$ws.Cells.Item(136, 9).Value = 2188.3125  # I136
$ws.Cells.Item(136, 11).Value = 6564.9375  # K136
$ws.Cells.Item(136, 13).Value = -4014.9375  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1802.6666  # H20
$ws.Cells.Item(20, 9).Value = 1802.6666  # I20
$ws.Cells.Item(20, 10).Value = 0  # J20
$ws.Cells.Item(20, 11).Value = 1802.6666  # K20
$ws.Cells.Item(20, 12).Value = 0  # L20
$ws.Cells.Item(20, 13).Value = -1555.6666  # M20
$ws.Cells.Item(20, 14).ClearContents()  # N20 removed
$ws.Cells.Item(86, 8).Value = 1880  # H86
$ws.Cells.Item(86, 9).Value = 1977.7778  # I86
$ws.Cells.Item(86, 10).Value = 1000  # J86
$ws.Cells.Item(86, 11).Value = 1977.7778  # K86
$ws.Cells.Item(86, 12).Value = 1000  # L86
$ws.Cells.Item(86, 13).Value = -854.7778000000001  # M86
$ws.Cells.Item(86, 14).Value = -3246  # N86
$ws.Cells.Item(89, 8).Value = 1880  # H89
$ws.Cells.Item(89, 9).Value = 1977.7778  # I89
$ws.Cells.Item(89, 10).Value = 1000  # J89
$ws.Cells.Item(89, 11).Value = 9888.889000000001  # K89
$ws.Cells.Item(89, 12).Value = 5000  # L89
$ws.Cells.Item(89, 13).Value = -4272.889000000001  # M89
$ws.Cells.Item(89, 14).Value = -16232  # N89
$ws.Cells.Item(105, 8).Value = 2715.6667  # H105
$ws.Cells.Item(105, 9).Value = 3000  # I105
$ws.Cells.Item(105, 10).Value = 2573.5  # J105
$ws.Cells.Item(105, 11).Value = 3000  # K105
$ws.Cells.Item(105, 12).Value = 2573.5  # L105
$ws.Cells.Item(105, 13).Value = -1253  # M105
$ws.Cells.Item(105, 14).Value = -6067.5  # N105

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2348.6875  # H16
$ws.Cells.Item(16, 9).Value = 2741.7778  # I16
$ws.Cells.Item(16, 11).Value = 2741.7778  # K16
$ws.Cells.Item(16, 13).Value = -2454.7778  # M16
$ws.Cells.Item(22, 8).Value = 575.5454999999999  # H22
$ws.Cells.Item(22, 10).Value = 652.6667  # J22
$ws.Cells.Item(22, 12).Value = 652.6667  # L22
$ws.Cells.Item(22, 14).Value = -1352.6667  # N22
$ws.Cells.Item(31, 8).Value = 1224.9231  # H31
$ws.Cells.Item(31, 9).Value = 1014.4737  # I31
$ws.Cells.Item(31, 11).Value = 1014.4737  # K31
$ws.Cells.Item(31, 13).Value = -719.4737  # M31
$ws.Cells.Item(34, 8).Value = 1224.9231  # H34
$ws.Cells.Item(34, 9).Value = 1014.4737  # I34
$ws.Cells.Item(34, 11).Value = 1014.4737  # K34
$ws.Cells.Item(34, 13).Value = -812.4737  # M34
$ws.Cells.Item(62, 8).Value = 4159.4  # H62
$ws.Cells.Item(62, 10).Value = 3999  # J62
$ws.Cells.Item(62, 12).Value = 3999  # L62
$ws.Cells.Item(62, 14).Value = -5247  # N62
$ws.Cells.Item(65, 8).Value = 4159.4  # H65
$ws.Cells.Item(65, 10).Value = 3999  # J65
$ws.Cells.Item(65, 12).Value = 19995  # L65
$ws.Cells.Item(65, 14).Value = -26235  # N65
$ws.Cells.Item(105, 8).Value = 2884.5  # H105
$ws.Cells.Item(105, 9).Value = 2347.5  # I105
$ws.Cells.Item(105, 10).Value = 3206.7  # J105
$ws.Cells.Item(105, 11).Value = 2347.5  # K105
$ws.Cells.Item(105, 12).Value = 3206.7  # L105
$ws.Cells.Item(105, 13).Value = -600.5  # M105
$ws.Cells.Item(105, 14).Value = -6700.7  # N105
$ws.Cells.Item(107, 8).Value = 821.619  # H107
$ws.Cells.Item(107, 9).Value = 915.6667  # I107
$ws.Cells.Item(107, 10).Value = 696.2222  # J107
$ws.Cells.Item(107, 11).Value = 915.6667  # K107
$ws.Cells.Item(107, 12).Value = 696.2222  # L107
$ws.Cells.Item(107, 13).Value = 1004.3333  # M107
$ws.Cells.Item(107, 14).Value = -4536.2222  # N107
$ws.Cells.Item(113, 8).Value = 2348.6875  # H113
$ws.Cells.Item(113, 9).Value = 2741.7778  # I113
$ws.Cells.Item(113, 11).Value = 2741.7778  # K113
$ws.Cells.Item(113, 13).Value = -571.7777999999998  # M113
$ws.Cells.Item(122, 8).Value = 2025.8334  # H122
$ws.Cells.Item(122, 9).Value = 2043.5714  # I122
$ws.Cells.Item(122, 11).Value = 6130.7142  # K122
$ws.Cells.Item(122, 13).Value = -3680.7142  # M122
$ws.Cells.Item(132, 8).Value = 3586.8  # H132
$ws.Cells.Item(132, 9).Value = 4167.2856  # I132
$ws.Cells.Item(132, 11).Value = 12501.8568  # K132
$ws.Cells.Item(132, 13).Value = -9971.856800000001  # M132
$ws.Cells.Item(134, 8).Value = 1025.1428  # H134
$ws.Cells.Item(134, 9).Value = 1046  # I134
$ws.Cells.Item(134, 10).Value = 900  # J134
$ws.Cells.Item(134, 11).Value = 3138  # K134
$ws.Cells.Item(134, 12).Value = 2700  # L134
$ws.Cells.Item(134, 13).Value = -603  # M134
$ws.Cells.Item(134, 14).Value = -7770  # N134
$ws.Cells.Item(141, 8).Value = 34991.08  # H141
$ws.Cells.Item(141, 9).Value = 0  # I141
$ws.Cells.Item(141, 10).Value = 34991.08  # J141
$ws.Cells.Item(141, 11).Value = 0  # K141
$ws.Cells.Item(141, 12).Value = 34991.08  # L141
$ws.Cells.Item(141, 13).ClearContents()  # M141 removed
$ws.Cells.Item(141, 14).Value = -45351.08  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 44.875  # H38
$ws.Cells.Item(38, 9).Value = 57.8  # I38
$ws.Cells.Item(38, 11).Value = 173.4  # K38
$ws.Cells.Item(38, 13).Value = 173.6  # M38
$ws.Cells.Item(63, 8).Value = 19483.334  # H63
$ws.Cells.Item(63, 9).Value = 19483.334  # I63
$ws.Cells.Item(63, 11).Value = 58450.00199999999  # K63
$ws.Cells.Item(63, 13).Value = -57701.00199999999  # M63
$ws.Cells.Item(66, 8).Value = 19483.334  # H66
$ws.Cells.Item(66, 9).Value = 19483.334  # I66
$ws.Cells.Item(66, 11).Value = 175350.006  # K66
$ws.Cells.Item(66, 13).Value = -171606.006  # M66
$ws.Cells.Item(69, 8).Value = 2228.4285  # H69
$ws.Cells.Item(69, 10).Value = 3000  # J69
$ws.Cells.Item(69, 12).Value = 9000  # L69
$ws.Cells.Item(69, 14).Value = -10622  # N69
$ws.Cells.Item(72, 8).Value = 2228.4285  # H72
$ws.Cells.Item(72, 10).Value = 3000  # J72
$ws.Cells.Item(72, 12).Value = 27000  # L72
$ws.Cells.Item(72, 14).Value = -35112  # N72
$ws.Cells.Item(75, 8).Value = 6619.6665  # H75
$ws.Cells.Item(75, 9).Value = 1251.5  # I75
$ws.Cells.Item(75, 10).Value = 9303.75  # J75
$ws.Cells.Item(75, 11).Value = 3754.5  # K75
$ws.Cells.Item(75, 12).Value = 27911.25  # L75
$ws.Cells.Item(75, 13).Value = -2756.5  # M75
$ws.Cells.Item(75, 14).Value = -29907.25  # N75
$ws.Cells.Item(78, 8).Value = 6619.6665  # H78
$ws.Cells.Item(78, 9).Value = 1251.5  # I78
$ws.Cells.Item(78, 10).Value = 9303.75  # J78
$ws.Cells.Item(78, 11).Value = 11263.5  # K78
$ws.Cells.Item(78, 12).Value = 83733.75  # L78
$ws.Cells.Item(78, 13).Value = -6271.5  # M78
$ws.Cells.Item(78, 14).Value = -93717.75  # N78
$ws.Cells.Item(92, 8).Value = 743.6667  # H92
$ws.Cells.Item(92, 10).Value = 799.1429000000001  # J92
$ws.Cells.Item(92, 12).Value = 2397.4287  # L92
$ws.Cells.Item(92, 14).Value = -4893.4287  # N92
$ws.Cells.Item(98, 8).Value = 1523.25  # H98
$ws.Cells.Item(98, 9).Value = 999  # I98
$ws.Cells.Item(98, 10).Value = 1598.1428  # J98
$ws.Cells.Item(98, 11).Value = 2997  # K98
$ws.Cells.Item(98, 12).Value = 4794.428400000001  # L98
$ws.Cells.Item(98, 13).Value = -1499  # M98
$ws.Cells.Item(98, 14).Value = -7790.428400000001  # N98
$ws.Cells.Item(114, 8).Value = 1960.8182  # H114
$ws.Cells.Item(114, 9).Value = 1064.875  # I114
$ws.Cells.Item(114, 10).Value = 4350  # J114
$ws.Cells.Item(114, 11).Value = 3194.625  # K114
$ws.Cells.Item(114, 12).Value = 13050  # L114
$ws.Cells.Item(114, 13).Value = 59.375  # M114
$ws.Cells.Item(114, 14).Value = -19558  # N114
$ws.Cells.Item(117, 8).Value = 25120.5  # H117
$ws.Cells.Item(117, 9).Value = 200  # I117
$ws.Cells.Item(117, 10).Value = 50041  # J117
$ws.Cells.Item(117, 11).Value = 600  # K117
$ws.Cells.Item(117, 12).Value = 150123  # L117
$ws.Cells.Item(117, 13).Value = 2842  # M117
$ws.Cells.Item(117, 14).Value = -157007  # N117
$ws.Cells.Item(131, 8).Value = 2043.35  # H131
$ws.Cells.Item(131, 10).Value = 2663.4614  # J131
$ws.Cells.Item(131, 12).Value = 7990.3842  # L131
$ws.Cells.Item(131, 14).Value = -18070.3842  # N131
$ws.Cells.Item(134, 8).Value = 7000  # H134
$ws.Cells.Item(134, 9).Value = 7000  # I134
$ws.Cells.Item(134, 10).Value = 0  # J134
$ws.Cells.Item(134, 11).Value = 21000  # K134
$ws.Cells.Item(134, 12).Value = 0  # L134
$ws.Cells.Item(134, 13).Value = -15930  # M134
$ws.Cells.Item(134, 14).ClearContents()  # N134 removed

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 50000  # H34
$ws.Cells.Item(34, 10).Value = 50000  # J34
$ws.Cells.Item(34, 12).Value = 50000  # L34
$ws.Cells.Item(34, 14).Value = -50536  # N34
$ws.Cells.Item(62, 8).Value = 90077  # H62
$ws.Cells.Item(62, 9).Value = 90077  # I62
$ws.Cells.Item(62, 11).Value = 90077  # K62
$ws.Cells.Item(62, 13).Value = -89391  # M62
$ws.Cells.Item(65, 8).Value = 90077  # H65
$ws.Cells.Item(65, 9).Value = 90077  # I65
$ws.Cells.Item(65, 11).Value = 270231  # K65
$ws.Cells.Item(65, 13).Value = -266799  # M65
$ws.Cells.Item(70, 8).Value = 3492.5  # H70
$ws.Cells.Item(70, 9).Value = 3492.5  # I70
$ws.Cells.Item(70, 11).Value = 3492.5  # K70
$ws.Cells.Item(70, 13).Value = -3222.5  # M70
$ws.Cells.Item(73, 8).Value = 3492.5  # H73
$ws.Cells.Item(73, 9).Value = 3492.5  # I73
$ws.Cells.Item(73, 11).Value = 3492.5  # K73
$ws.Cells.Item(73, 13).Value = -2556.5  # M73
$ws.Cells.Item(76, 8).Value = 50000  # H76
$ws.Cells.Item(76, 10).Value = 50000  # J76
$ws.Cells.Item(76, 12).Value = 50000  # L76
$ws.Cells.Item(76, 14).Value = -50630  # N76
$ws.Cells.Item(79, 8).Value = 50000  # H79
$ws.Cells.Item(79, 10).Value = 50000  # J79
$ws.Cells.Item(79, 12).Value = 50000  # L79
$ws.Cells.Item(79, 14).Value = -52184  # N79
$ws.Cells.Item(97, 8).Value = 665  # H97
$ws.Cells.Item(97, 9).Value = 665  # I97
$ws.Cells.Item(97, 11).Value = 665  # K97
$ws.Cells.Item(97, 13).Value = -169  # M97
$ws.Cells.Item(113, 8).Value = 961.1  # H113
$ws.Cells.Item(113, 9).Value = 961.1  # I113
$ws.Cells.Item(113, 11).Value = 961.1  # K113
$ws.Cells.Item(113, 13).Value = 1208.9  # M113
$ws.Cells.Item(126, 8).Value = 3577.4  # H126
$ws.Cells.Item(126, 9).Value = 3577.4  # I126
$ws.Cells.Item(126, 11).Value = 10732.2  # K126
$ws.Cells.Item(126, 13).Value = -8262.200000000001  # M126

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 573.9231  # H16
$ws.Cells.Item(16, 9).Value = 588.5  # I16
$ws.Cells.Item(16, 11).Value = 588.5  # K16
$ws.Cells.Item(16, 13).Value = -418.5  # M16
$ws.Cells.Item(63, 8).Value = 84077  # H63
$ws.Cells.Item(63, 9).Value = 84077  # I63
$ws.Cells.Item(63, 11).Value = 84077  # K63
$ws.Cells.Item(63, 13).Value = -83328  # M63
$ws.Cells.Item(66, 8).Value = 84077  # H66
$ws.Cells.Item(66, 9).Value = 84077  # I66
$ws.Cells.Item(66, 11).Value = 252231  # K66
$ws.Cells.Item(66, 13).Value = -248487  # M66
$ws.Cells.Item(68, 8).Value = 2458.1667  # H68
$ws.Cells.Item(68, 10).Value = 1899  # J68
$ws.Cells.Item(68, 12).Value = 1899  # L68
$ws.Cells.Item(68, 14).Value = -3397  # N68
$ws.Cells.Item(71, 8).Value = 2458.1667  # H71
$ws.Cells.Item(71, 10).Value = 1899  # J71
$ws.Cells.Item(71, 12).Value = 9495  # L71
$ws.Cells.Item(71, 14).Value = -16983  # N71
$ws.Cells.Item(93, 8).Value = 3793.6  # H93
$ws.Cells.Item(93, 9).Value = 3793.6  # I93
$ws.Cells.Item(93, 11).Value = 3793.6  # K93
$ws.Cells.Item(93, 13).Value = -2545.6  # M93

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(57, 8).Value = 115000  # H57
$ws.Cells.Item(57, 9).Value = 115000  # I57
$ws.Cells.Item(57, 11).Value = 115000  # K57
$ws.Cells.Item(57, 13).Value = -114246  # M57
$ws.Cells.Item(62, 8).Value = 23056.857  # H62
$ws.Cells.Item(62, 10).Value = 10278.8  # J62
$ws.Cells.Item(62, 12).Value = 10278.8  # L62
$ws.Cells.Item(62, 14).Value = -11526.8  # N62
$ws.Cells.Item(65, 8).Value = 23056.857  # H65
$ws.Cells.Item(65, 10).Value = 10278.8  # J65
$ws.Cells.Item(65, 12).Value = 51394  # L65
$ws.Cells.Item(65, 14).Value = -57634  # N65
$ws.Cells.Item(81, 8).Value = 1429060.8  # H81
$ws.Cells.Item(81, 9).Value = 531  # I81
$ws.Cells.Item(81, 10).Value = 3333767  # J81
$ws.Cells.Item(81, 11).Value = 1062  # K81
$ws.Cells.Item(81, 12).Value = 6667534  # L81
$ws.Cells.Item(81, 13).Value = -1  # M81
$ws.Cells.Item(81, 14).Value = -6669656  # N81
$ws.Cells.Item(84, 8).Value = 1429060.8  # H84
$ws.Cells.Item(84, 9).Value = 531  # I84
$ws.Cells.Item(84, 10).Value = 3333767  # J84
$ws.Cells.Item(84, 11).Value = 5310  # K84
$ws.Cells.Item(84, 12).Value = 33337670  # L84
$ws.Cells.Item(84, 13).Value = -6  # M84
$ws.Cells.Item(84, 14).Value = -33348278  # N84
$ws.Cells.Item(107, 8).Value = 601.9167  # H107
$ws.Cells.Item(107, 9).Value = 491.66666  # I107
$ws.Cells.Item(107, 11).Value = 1474.99998  # K107
$ws.Cells.Item(107, 13).Value = 445.0000199999999  # M107
$ws.Cells.Item(110, 8).Value = 64400  # H110
$ws.Cells.Item(110, 10).Value = 64400  # J110
$ws.Cells.Item(110, 12).Value = 64400  # L110
$ws.Cells.Item(110, 14).Value = -72580  # N110
$ws.Cells.Item(120, 8).Value = 20000  # H120
$ws.Cells.Item(120, 10).Value = 20000  # J120
$ws.Cells.Item(120, 12).Value = 20000  # L120
$ws.Cells.Item(120, 14).Value = -29676  # N120
$ws.Cells.Item(126, 8).Value = 2820  # H126
$ws.Cells.Item(126, 9).Value = 2820  # I126
$ws.Cells.Item(126, 11).Value = 8460  # K126
$ws.Cells.Item(126, 13).Value = -5990  # M126
$ws.Cells.Item(132, 8).Value = 2970  # H132
$ws.Cells.Item(132, 9).Value = 4166.3335  # I132
$ws.Cells.Item(132, 11).Value = 12499.0005  # K132
$ws.Cells.Item(132, 13).Value = -9969.000499999998  # M132
$ws.Cells.Item(136, 8).Value = 1894.12  # H136
$ws.Cells.Item(136, 9).Value = 1863.9546  # I136
$ws.Cells.Item(136, 11).Value = 5591.8638  # K136
$ws.Cells.Item(136, 13).Value = -3041.8638  # M136
